$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(45986, "FRAMBUESA", "COOPERATIVA HORTOFRUTICOLA CARTAYA", "Z1", 351, 3294, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "COOPERATIVA HORTOFRUTICOLA CARTAYA", "Z2", 1203, 32339.7, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "COOPERATIVA HORTOFRUTICOLA CARTAYA", "Z3", 0, 42393, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "COOPERATIVA HORTOFRUTICOLA CARTAYA", "Z4", 597, 27924, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "FRESLUCENA, S.A.", "Z1", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "FRESLUCENA, S.A.", "Z2", 0, 667.5, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "FRESLUCENA, S.A.", "Z3", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "FRESLUCENA, S.A.", "Z4", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "ONUBAFRUIT", "Z1", 351, 3402, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "ONUBAFRUIT", "Z2", 1731, 52385.7, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "ONUBAFRUIT", "Z3", 387, 56688, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "ONUBAFRUIT", "Z4", 1245, 30147, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.A.T CONDADO DE HUELVA", "Z1", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.A.T CONDADO DE HUELVA", "Z2", 0, 1074, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.A.T CONDADO DE HUELVA", "Z3", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.A.T CONDADO DE HUELVA", "Z4", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. COSTA DE HUELVA", "Z1", 0, 108, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. COSTA DE HUELVA", "Z2", 528, 8740.5, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. COSTA DE HUELVA", "Z3", 0, 3114, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. COSTA DE HUELVA", "Z4", 648, 2223, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. NUESTRA SRA. LA BELLA", "Z1", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. NUESTRA SRA. LA BELLA", "Z2", 0, 9564, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. NUESTRA SRA. LA BELLA", "Z3", 387, 11181, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "S.C.A. NUESTRA SRA. LA BELLA", "Z4", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "VARIOS", "Z1", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "VARIOS", "Z2", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "VARIOS", "Z3", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
    ,@(45986, "FRAMBUESA", "VARIOS", "Z4", 0, 0, "12 DATOS SEMANA ONUBAFRUIT 13112025 A 19112025.xlsx")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value2 = $item[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
}

Write-Host "Done adding rows"
$usedAddr = $ws.UsedRange.Address()
Write-Host $usedAddr
